$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.853.85"
$ws.Range("E2").Value = "  +3.24%  "

$ws.Range("D3").Value = "2.535.34"
$ws.Range("E3").Value = "  +2.79%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "581.71"
$ws.Range("E5").Value = "  +1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.10"
$ws.Range("E6").Value = "  +3.66%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +1.13%  "

$ws.Range("D9").Value = "2.538.46"
$ws.Range("E9").Value = "  +2.91%  "

$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").Value = "29.32"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("D16").Value = "2.999.03"
$ws.Range("E16").Value = "  +2.84%  "

$ws.Range("D17").Value = "64.326.48"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").Value = "2.539.83"
$ws.Range("E18").Value = "  +2.95%  "

$ws.Range("D19").Value = "8.01"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").Value = "11.02"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").Value = "4.29"
$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("D22").Value = "330.02"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").Value = "2.25"
$ws.Range("E23").Value = "  +1.80%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "10.17"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").Value = "65.83"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").Value = "648.56"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("E28").Value = "  +6.68%  "

$ws.Range("D29").Value = "2.644.06"
$ws.Range("E29").Value = "  +2.03%  "

$ws.Range("E30").Value = "  +4.57%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("D34").Value = "0.138"
$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  +1.67%  "

$ws.Range("E37").Value = "  +2.34%  "

$ws.Range("D38").Value = "5.62"
$ws.Range("E38").Value = "  +4.49%  "

$ws.Range("D39").Value = "155.92"
$ws.Range("E39").Value = "  +2.45%  "

$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "19.02"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.373"
$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").Value = "1.82"
$ws.Range("E43").Value = "  +4.77%  "

$ws.Range("D44").Value = "163.22"
$ws.Range("E44").Value = "  +6.70%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "0.0₆0303"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").Value = "15.64"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").Value = "3.67"
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("D49").Value = "21.66"
$ws.Range("E49").Value = "  +5.78%  "

$ws.Range("D50").Value = "0.625"
$ws.Range("E50").Value = "  +2.67%  "

$ws.Range("E51").Value = "  +1.41%  "
